$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the last existing data row (357) down through the new rows
$ws.Range("A357:D357").Copy()
$ws.Range("A358:D366").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A358").Value = 44432
$ws.Range("B358").Value = 41
$ws.Range("C358").Value = 190
$ws.Range("D358").Value = 100.5221863046457

$ws.Range("A359").Value = 44433
$ws.Range("B359").Value = 22
$ws.Range("C359").Value = 199
$ws.Range("D359").Value = 105.2837635506552

$ws.Range("A360").Value = 44434
$ws.Range("B360").Value = 7
$ws.Range("C360").Value = 179
$ws.Range("D360").Value = 94.70248078174517

$ws.Range("A361").Value = 44435
$ws.Range("B361").Value = 36
$ws.Range("C361").Value = 186
$ws.Range("D361").Value = 98.40592975086369

$ws.Range("A362").Value = 44436
$ws.Range("B362").Value = 28
$ws.Range("C362").Value = 173
$ws.Range("D362").Value = 91.52809595107215

$ws.Range("A363").Value = 44437
$ws.Range("B363").Value = 23
$ws.Range("C363").Value = 188
$ws.Range("D363").Value = 99.46405802775472

$ws.Range("A364").Value = 44438
$ws.Range("B364").Value = 26
$ws.Range("C364").Value = 183
$ws.Range("D364").Value = 96.81873733552719

$ws.Range("A365").Value = 44439
$ws.Range("B365").Value = 27
$ws.Range("C365").Value = 169
$ws.Range("D365").Value = 89.41183939729014

$ws.Range("A366").Value = 44440
$ws.Range("B366").Value = 14
$ws.Range("C366").Value = 161
$ws.Range("D366").Value = 85.17932628972611
